# Weekly update: insert a new price record for "Acelga" (Macroferia Regional
# de Talca) as row 271, pushing the existing rows 271-388 down to 272-389.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 271 (shifts rows 271..388 down to 272..389).
$ws.Rows.Item(271).Insert()

# Populate the newly inserted row with the new week's data.
$ws.Cells.Item(271, 1).Value  = 5
$ws.Cells.Item(271, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(271, 3).Value  = "Maule"
$ws.Cells.Item(271, 4).Value  = 44992
$ws.Cells.Item(271, 5).Value  = 7
$ws.Cells.Item(271, 6).Value  = 100112009
$ws.Cells.Item(271, 7).Value  = "Acelga"
$ws.Cells.Item(271, 8).Value  = "Sin especificar"
$ws.Cells.Item(271, 9).Value  = "Primera"
$ws.Cells.Item(271, 10).Value = 500
$ws.Cells.Item(271, 11).Value = 3000
$ws.Cells.Item(271, 12).Value = 3000
$ws.Cells.Item(271, 13).Value = 3000
$ws.Cells.Item(271, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(271, 15).Value = "Región del Maule"
$ws.Cells.Item(271, 16).Value = 750
$ws.Cells.Item(271, 17).Value = 4
$ws.Cells.Item(271, 18).Value = "Hortaliza"
